# Delete rows 7 and 8 (the "SMS Respon" / "SMS Test2" row and the
# "same" / "Clock" row) so that the remaining rows shift up, matching
# the AgentSettingsWMC test data expectations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("7:8").Delete()
